# Apply "想去人数" (F column) updates to the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row number -> new value }
$updates = @{
    "展览" = @{
        5  = 119
        10 = 16327
        12 = 202
        13 = 1038
        14 = 6360
        15 = 637
        17 = 80
        19 = 125
        21 = 34
        24 = 33
        30 = 5049
        32 = 11323
        36 = 206
        37 = 3836
        39 = 74
    }
    "全部类型" = @{
        5  = 119
        10 = 16327
        12 = 202
        13 = 1038
        14 = 6360
        15 = 637
        17 = 80
        19 = 125
        21 = 34
        24 = 33
        30 = 5049
        33 = 11323
        37 = 206
        38 = 3836
        40 = 74
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
